$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new "Variance" column before the existing "Integral" column (old C).
# This shifts old C->D (Integral) and old D->E (Time); formulas/refs shift too.
$ws.Columns("C").Insert()

# Header + data for the new Variance column (= STD^2)
$ws.Range("C1").Value = "Variance"
$ws.Range("C2").Formula = "=B2^2"
$ws.Range("C3:C11").Formula = "=B3^2"

# Average of the new Variance column
$ws.Range("C13").Formula = "=AVERAGE(C2:C11)"

# Remove the old STD formulas for the raw STD column (B) and new Variance column (C);
# only Integral (D) and Time (E) keep a STD row now.
$ws.Range("B14").ClearContents()
$ws.Range("C14").ClearContents()

# New RMS row: RMS = sqrt(average variance)
$ws.Range("A15").Value = "RMS"
$ws.Range("B15").Formula = "=SQRT(C13)"

# Column C width (Excel auto "best fit" sizing, ~11.16 chars wide)
$ws.Columns("C").ColumnWidth = 10.4

# Match the saved selection/active cell from the source workbook
$ws.Range("B16").Select() | Out-Null
